# Updates currentAveragePrice / profit figures across multiple sheets
# per the scheduled pricing-refresh run (H/I/J/K/L/M/N columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 240.75
$ws.Range("I33").Value = 248.4375
$ws.Range("K33").Value = 248.4375
$ws.Range("M33").Value = -19.4375
$ws.Range("H58").Value = 2455.84
$ws.Range("J58").Value = 3222.875
$ws.Range("L58").Value = 9668.625
$ws.Range("N58").Value = -9968.625
$ws.Range("H86").Value = 2998
$ws.Range("I86").Value = 1996.8334
$ws.Range("K86").Value = 1996.8334
$ws.Range("M86").Value = -873.8334
$ws.Range("H89").Value = 2998
$ws.Range("I89").Value = 1996.8334
$ws.Range("K89").Value = 9984.166999999999
$ws.Range("M89").Value = -4368.166999999999
$ws.Range("H92").Value = 1540.9375
$ws.Range("I92").Value = 977
$ws.Range("K92").Value = 977
$ws.Range("M92").Value = 271
$ws.Range("H129").Value = 1627.1111
$ws.Range("I129").Value = 1027.5
$ws.Range("K129").Value = 3082.5
$ws.Range("M129").Value = 1917.5
$ws.Range("H135").Value = 1815.1904
$ws.Range("I135").Value = 1058.1428
$ws.Range("J135").Value = 3329.2856
$ws.Range("K135").Value = 9523.2852
$ws.Range("L135").Value = 29963.5704
$ws.Range("M135").Value = -6988.2852
$ws.Range("N135").Value = -35033.5704
$ws.Range("H137").Value = 3236.764
$ws.Range("I137").Value = 2167.7754
$ws.Range("J137").Value = 5514.174
$ws.Range("K137").Value = 6503.3262
$ws.Range("L137").Value = 16542.522
$ws.Range("M137").Value = -3953.3262
$ws.Range("N137").Value = -21642.522
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21560.219
$ws.Range("I32").Value = 21265.938
$ws.Range("J32").Value = 22606.555
$ws.Range("K32").Value = 21265.938
$ws.Range("L32").Value = 22606.555
$ws.Range("M32").Value = -20978.938
$ws.Range("N32").Value = -23180.555
$ws.Range("H74").Value = 289432.16
$ws.Range("J74").Value = 8592.299999999999
$ws.Range("L74").Value = 8592.299999999999
$ws.Range("N74").Value = -10340.3
$ws.Range("H77").Value = 289432.16
$ws.Range("J77").Value = 8592.299999999999
$ws.Range("L77").Value = 42961.5
$ws.Range("N77").Value = -51697.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 3610.2
$ws.Range("I22").Value = 2262.75
$ws.Range("K22").Value = 2262.75
$ws.Range("M22").Value = -2089.75
$ws.Range("H81").Value = 41602.31
$ws.Range("J81").Value = 41602.31
$ws.Range("L81").Value = 41602.31
$ws.Range("N81").Value = -43724.31
$ws.Range("H84").Value = 41602.31
$ws.Range("J84").Value = 41602.31
$ws.Range("L84").Value = 124806.93
$ws.Range("N84").Value = -135414.93
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 50556
$ws.Range("J20").Value = 50556
$ws.Range("L20").Value = 50556
$ws.Range("N20").Value = -51028
$ws.Range("H30").Value = 50556
$ws.Range("J30").Value = 50556
$ws.Range("L30").Value = 50556
$ws.Range("N30").Value = -50738
$ws.Range("H31").Value = 31255058
$ws.Range("I31").Value = 52634130
$ws.Range("K31").Value = 52634130
$ws.Range("M31").Value = -52633835
$ws.Range("H34").Value = 31255058
$ws.Range("I34").Value = 52634130
$ws.Range("K34").Value = 52634130
$ws.Range("M34").Value = -52633928
$ws.Range("H58").Value = 6397.8335
$ws.Range("I58").Value = 3332.7144
$ws.Range("K58").Value = 3332.7144
$ws.Range("M58").Value = -3129.7144
$ws.Range("H94").Value = 2849.8
$ws.Range("J94").Value = 1155.2142
$ws.Range("L94").Value = 1155.2142
$ws.Range("N94").Value = -2057.2142
$ws.Range("H128").Value = 50556
$ws.Range("J128").Value = 50556
$ws.Range("L128").Value = 50556
$ws.Range("N128").Value = -60516
$ws.Range("H132").Value = 14473.981
$ws.Range("I132").Value = 3284.8462
$ws.Range("K132").Value = 9854.5386
$ws.Range("M132").Value = -7324.5386
$ws.Range("H136").Value = 6397.8335
$ws.Range("I136").Value = 3332.7144
$ws.Range("K136").Value = 9998.143199999999
$ws.Range("M136").Value = -7448.143199999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8043.8
$ws.Range("I3").Value = 6204.385
$ws.Range("K3").Value = 18613.155
$ws.Range("M3").Value = -18501.155
$ws.Range("H5").Value = 1130.7142
$ws.Range("I5").Value = 449.0909
$ws.Range("J5").Value = 3630
$ws.Range("K5").Value = 1347.2727
$ws.Range("L5").Value = 10890
$ws.Range("M5").Value = -1235.2727
$ws.Range("N5").Value = -11114
$ws.Range("H135").Value = 1130.7142
$ws.Range("I135").Value = 449.0909
$ws.Range("J135").Value = 3630
$ws.Range("K135").Value = 4041.8181
$ws.Range("L135").Value = 32670
$ws.Range("M135").Value = -1506.8181
$ws.Range("N135").Value = -37740
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 951
$ws.Range("J12").Value = 2500
$ws.Range("L12").Value = 2500
$ws.Range("N12").Value = -2780
$ws.Range("H40").Value = 25832.666
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 25832.666
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 25832.666
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -26134.666
$ws.Range("H43").Value = 1410.8889
$ws.Range("I43").Value = 1410.8889
$ws.Range("K43").Value = 1410.8889
$ws.Range("M43").Value = -1259.8889
$ws.Range("H122").Value = 3555
$ws.Range("J122").Value = 4928
$ws.Range("L122").Value = 14784
$ws.Range("N122").Value = -19684
$ws.Range("H126").Value = 2964.7646
$ws.Range("I126").Value = 2233.6
$ws.Range("K126").Value = 6700.799999999999
$ws.Range("M126").Value = -4230.799999999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3707.1667
$ws.Range("I7").Value = 2469.5715
$ws.Range("K7").Value = 2469.5715
$ws.Range("M7").Value = -2357.5715
$ws.Range("H40").Value = 26317994
$ws.Range("I40").Value = 33335372
$ws.Range("J40").Value = 2823
$ws.Range("K40").Value = 33335372
$ws.Range("L40").Value = 2823
$ws.Range("M40").Value = -33335236
$ws.Range("N40").Value = -3095
$ws.Range("H46").Value = 6960.483
$ws.Range("J46").Value = 9026.619000000001
$ws.Range("L46").Value = 9026.619000000001
$ws.Range("N46").Value = -9402.619000000001
$ws.Range("H126").Value = 3707.1667
$ws.Range("I126").Value = 2469.5715
$ws.Range("K126").Value = 7408.7145
$ws.Range("M126").Value = -4938.7145
$ws.Range("H132").Value = 5322.26
$ws.Range("I132").Value = 4596.2188
$ws.Range("J132").Value = 6613
$ws.Range("K132").Value = 13788.6564
$ws.Range("L132").Value = 19839
$ws.Range("M132").Value = -11258.6564
$ws.Range("N132").Value = -24899
$ws.Range("H136").Value = 1555693.2
$ws.Range("I136").Value = 2781219.2
$ws.Range("K136").Value = 8343657.600000001
$ws.Range("M136").Value = -8341107.600000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 2100
$ws.Range("J12").Value = 2100
$ws.Range("L12").Value = 2100
$ws.Range("N12").Value = -2384
$ws.Range("H107").Value = 1433.3334
$ws.Range("I107").Value = 1433.3334
$ws.Range("K107").Value = 4300.0002
$ws.Range("M107").Value = -2380.0002
$ws.Range("H126").Value = 6213.303
$ws.Range("I126").Value = 3509.7693
$ws.Range("K126").Value = 10529.3079
$ws.Range("M126").Value = -8059.3079
